# Add new rows of platform/url data to the "valid_csv" sheet, matching the
# rows appended to the sample upload fixture (YouTube/TikTok links used to
# exercise the new metadata-fetcher flow).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("valid_csv")

$newRows = @(
    @("youtube", "https://youtu.be/fDJ2NrRRoB8?si=8Jpe-jjjUCR0Vd2m"),
    @("youtube", "https://www.youtube.com/watch?v=jcb6Lygf3LQ"),
    @("youtube", "https://www.youtube.com/watch?v=W69ZXgHm65A&list=RDW69ZXgHm65A&start_radio=1"),
    @("tiktok",  "https://www.tiktok.com/@paramountpics/video/7551460293851811103"),
    @("tiktok",  "https://www.tiktok.com/@aleko.so/video/7556646272849956107"),
    @("youtube", "https://www.youtube.com/shorts/5LekKnvD83E")
)

$row = 6
foreach ($pair in $newRows) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Match the window/view tweaks captured in the diff.
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("A15").Select()
